$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.908.72'
$ws.Range('E2').Value = '  +2.80%  '
$ws.Range('D3').Value = '3.448.56'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.24'
$ws.Range('E5').Value = '  +4.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '187.13'
$ws.Range('E6').Value = '  +6.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.444.20'
$ws.Range('E8').Value = '  +2.40%  '
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E10').Value = '  -1.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.646'
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.44'
$ws.Range('E12').Value = '  +4.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000277'
$ws.Range('E13').Value = '  -1.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.42'
$ws.Range('E14').Value = '  +2.86%  '
$ws.Range('D15').Value = '3.996.03'
$ws.Range('E15').Value = '  +2.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.71'
$ws.Range('E16').Value = '  +2.43%  '
$ws.Range('D17').Value = '3.444.27'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('D18').Value = '66.845.71'
$ws.Range('E18').Value = '  +2.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.10'
$ws.Range('E19').Value = '  +1.81%  '
$ws.Range('E20').Value = '  -2.10%  '
$ws.Range('E21').Value = '  +2.69%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '484.91'
$ws.Range('E22').Value = '  +6.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.31'
$ws.Range('E23').Value = '  +8.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.85'
$ws.Range('E24').Value = '  +22.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.35'
$ws.Range('E25').Value = '  +6.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.38'
$ws.Range('E26').Value = '  +2.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.96'
$ws.Range('E27').Value = '  +2.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.00'
$ws.Range('E28').Value = '  +1.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.06'
$ws.Range('E29').Value = '  +4.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.25'
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.25'
$ws.Range('E31').Value = '  +10.20%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.74'
$ws.Range('E32').Value = '  +2.37%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '596.32'
$ws.Range('E33').Value = '  +3.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.45'
$ws.Range('E34').Value = '  +0.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.112'
$ws.Range('E35').Value = '  +3.57%  '
$ws.Range('E36').Value = '  +5.93%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.72'
$ws.Range('E38').Value = '  +2.68%  '
$ws.Range('E39').Value = '  +3.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.54'
$ws.Range('E40').Value = '  -2.67%  '
$ws.Range('D41').Value = '3.258.54'
$ws.Range('D42').Value = '0.0₃0752'
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('E43').Value = '  +4.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0431'
$ws.Range('E44').Value = '  +2.95%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.53'
$ws.Range('E45').Value = '  +2.86%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.79'
$ws.Range('E46').Value = '  +22.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.23'
$ws.Range('E47').Value = '  +0.74%  '
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.34'
$ws.Range('E49').Value = '  +14.82%  '
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.71'
